$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.857.55'
$ws.Range("E2").Value = '  +4.13%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.276.08'
$ws.Range("E3").Value = '  +2.20%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.30'
$ws.Range("E5").Value = '  +4.26%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '92.69'
$ws.Range("E6").Value = '  +4.95%  '
$ws.Range("E7").Value = '  +3.72%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.486'
$ws.Range("E9").Value = '  +3.31%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.71'
$ws.Range("E10").Value = '  +6.16%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.76'
$ws.Range("E11").Value = '  +5.71%  '
$ws.Range("E12").Value = '  +2.01%  '
$ws.Range("E13").Value = '  +1.37%  '
$ws.Range("E14").Value = '  +3.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.628.05'
$ws.Range("E15").Value = '  +3.27%  '
$ws.Range("E16").Value = '  +3.32%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.257.40'
$ws.Range("E17").Value = '  +2.33%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.765'
$ws.Range("E18").Value = '  +3.53%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '41.797.47'
$ws.Range("E19").Value = '  +4.16%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.34'
$ws.Range("E20").Value = '  +9.10%  '
$ws.Range("E21").Value = '  +2.13%  '
$ws.Range("E22").Value = '  +2.53%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.16'
$ws.Range("E23").Value = '  +2.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '242.74'
$ws.Range("E24").Value = '  +2.72%  '
$ws.Range("E25").Value = '  +5.22%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.93'
$ws.Range("E27").Value = '  +5.19%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '24.35'
$ws.Range("E28").Value = '  +4.39%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.63'
$ws.Range("E29").Value = '  +2.93%  '
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.36'
$ws.Range("E31").Value = '  +7.29%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '158.92'
$ws.Range("E32").Value = '  -0.10%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  +0.04%  '
$ws.Range("E34").Value = '  +4.19%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0751'
$ws.Range("E35").Value = '  +4.64%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.05'
$ws.Range("E36").Value = '  +1.01%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '17.14'
$ws.Range("E37").Value = '  +8.87%  '
$ws.Range("E39").Value = '  +2.61%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.105'
$ws.Range("E40").Value = '  +5.11%  '
$ws.Range("E41").Value = '  +3.05%  '
$ws.Range("E42").Value = '  +4.42%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.074.63'
$ws.Range("E43").Value = '  -0.76%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.46'
$ws.Range("E44").Value = '  +0.13%  '
$ws.Range("E45").Value = '  +3.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.37'
$ws.Range("E46").Value = '  +2.68%  '
$ws.Range("E47").Value = '  +5.31%  '
$ws.Range("E48").Value = '  +7.87%  '
$ws.Range("E49").Value = '  +3.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '73.08'
$ws.Range("E50").Value = '  +7.19%  '
$ws.Range("E51").Value = '  +3.44%  '
